$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.837.17"
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").Value = "3.495.59"
$ws.Range("E3").Value = "  -2.09%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.41"
$ws.Range("E5").Value = "  -0.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "197.49"
$ws.Range("E6").Value = "  +5.79%  "

$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -3.47%  "

$ws.Range("E10").Value = "  +0.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.08"
$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000302"
$ws.Range("E12").Value = "  -2.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.55"
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").Value = "4.054.78"
$ws.Range("E14").Value = "  -1.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "592.92"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "69.906.75"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.00"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D19").Value = "3.498.58"
$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.985"
$ws.Range("E21").Value = "  -1.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.01"
$ws.Range("E22").Value = "  +3.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "103.86"
$ws.Range("E23").Value = "  +9.87%  "

$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.00"
$ws.Range("E25").Value = "  +2.28%  "

$ws.Range("E26").Value = "  +4.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("E28").Value = "  +3.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.43"
$ws.Range("E29").Value = "  +3.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.48"
$ws.Range("E30").Value = "  +20.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.27"
$ws.Range("E31").Value = "  +2.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.71"
$ws.Range("E32").Value = "  +3.58%  "

$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.54"
$ws.Range("E34").Value = "  -0.22%  "

$ws.Range("D35").Value = "3.703.20"
$ws.Range("E35").Value = "  +4.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("D37").Value = "0.0₃0802"
$ws.Range("E37").Value = "  +1.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "511.28"
$ws.Range("E38").Value = "  -3.85%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.98"
$ws.Range("E39").Value = "  -7.24%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.390"
$ws.Range("E40").Value = "  -3.89%  "

$ws.Range("E41").Value = "  -2.53%  "

$ws.Range("E42").Value = "  -1.04%  "

$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.83"
$ws.Range("E45").Value = "  -3.42%  "

$ws.Range("E46").Value = "  -1.47%  "

$ws.Range("E47").Value = "  -4.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.74"
$ws.Range("E49").Value = "  -5.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.44"
$ws.Range("E50").Value = "  -3.13%  "

$ws.Range("E51").Value = "  -3.00%  "
